$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Philadelphia, Pennsylvania (row 12)
$ws.Range("B12").Value = 30.75489282385834
$ws.Range("C12").Value = 25.98586704134461
$ws.Range("D12").Value = -0.1550655958980979

# Baltimore, Maryland (row 13)
$ws.Range("B13").Value = 57.02947845804988
$ws.Range("C13").Value = 41.22492887372787
$ws.Range("D13").Value = -0.2771294778006365

# Atlanta, Georgia (row 14)
$ws.Range("B14").Value = 31.47632311977716
$ws.Range("C14").Value = 26.38881114407491
$ws.Range("D14").Value = -0.1616298052457616
